# WebApi - Added new competitions and special offer
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Re-type the match dates on the existing Croatian/English fixtures so
#     they use a leading-zero day (e.g. "4.11." -> "04.11."); this mirrors
#     the author's retyped MatchDateTime values and also flips the
#     SpecialOffer flag on several of them. ---
$ws.Range("K2").Value = "04.11. 18:00"
$ws.Range("L2").Value = $true

$ws.Range("K3").Value = "04.11. 18:00"

$ws.Range("K4").Value = "05.11. 17:30"
$ws.Range("L4").Value = $true

$ws.Range("K6").Value = "06.11. 17:10"
$ws.Range("L6").Value = $true

$ws.Range("K7").Value = "05.11. 16:00"
$ws.Range("L7").Value = $true

$ws.Range("K8").Value = "05.11. 16:00"
$ws.Range("L8").Value = $true

$ws.Range("K9").Value = "05.11. 16:00"
$ws.Range("L9").Value = $true

$ws.Range("K10").Value = "05.11. 16:00"
$ws.Range("L10").Value = $true

$ws.Range("K11").Value = "05.11. 16:00"

# --- Helper to fill one betting-data row with the shared layout used by
#     the sheet (A=MatchId, B=Competition, C=HomeTeam, D=AwayTeam,
#     E..J=odds, K=MatchDateTime, L=SpecialOffer). Positional params only -
#     this COM host does not bind named (-Param value) arguments. ---
function Set-BettingRow {
    param(
        [int]$Row,
        [int]$MatchId,
        [string]$Competition,
        [string]$HomeTeam,
        [string]$AwayTeam,
        [double]$HomeWin,
        [double]$Draw,
        [double]$AwayWin,
        [double]$HomeOrDraw,
        [double]$AwayOrDraw,
        [double]$HomeOrAway,
        [string]$MatchDateTime,
        [bool]$SpecialOffer
    )

    $ws.Cells.Item($Row, 1).Value = $MatchId
    $ws.Cells.Item($Row, 2).Value = $Competition
    $ws.Cells.Item($Row, 3).Value = $HomeTeam
    $ws.Cells.Item($Row, 4).Value = $AwayTeam

    $ws.Cells.Item($Row, 5).NumberFormat = $ws.Cells.Item(2, 5).NumberFormat
    $ws.Cells.Item($Row, 5).Value = $HomeWin
    $ws.Cells.Item($Row, 6).NumberFormat = $ws.Cells.Item(2, 6).NumberFormat
    $ws.Cells.Item($Row, 6).Value = $Draw
    $ws.Cells.Item($Row, 7).NumberFormat = $ws.Cells.Item(2, 7).NumberFormat
    $ws.Cells.Item($Row, 7).Value = $AwayWin
    $ws.Cells.Item($Row, 8).NumberFormat = $ws.Cells.Item(2, 8).NumberFormat
    $ws.Cells.Item($Row, 8).Value = $HomeOrDraw
    $ws.Cells.Item($Row, 9).NumberFormat = $ws.Cells.Item(2, 9).NumberFormat
    $ws.Cells.Item($Row, 9).Value = $AwayOrDraw
    $ws.Cells.Item($Row, 10).NumberFormat = $ws.Cells.Item(2, 10).NumberFormat
    $ws.Cells.Item($Row, 10).Value = $HomeOrAway

    $ws.Cells.Item($Row, 11).Value = $MatchDateTime
    $ws.Cells.Item($Row, 12).Value = $SpecialOffer
}

# --- New competition: Spanjolska 1 (Spain) ---
Set-BettingRow 12 11 "Španjolska 1" "Espanyol" "Villarreal" `
    2.85 3.3 2.6 1.55 1.45 1.35 "09.11. 20:00" $false

Set-BettingRow 13 12 "Španjolska 1" "Real Madrid" "Cadiz" `
    1.15 8.2 17 1.02 5.5 1.1 "10.11. 20:30" $false

Set-BettingRow 14 13 "Španjolska 1" "Atl. Madrid" "Elche" `
    1.3 5.6 10 1.05 3.6 1.15 "10.11. 20:00" $true

Set-BettingRow 15 14 "Španjolska 1" "Valencia" "Betis" `
    2.2 3.4 3.4 1.35 1.7 1.35 "10.11. 20:30" $false

# --- New competition: Italija 1 (Italy) ---
Set-BettingRow 16 15 "Italija 1" "Fiorentina" "Salernitana" `
    1.55 4.3 6 1.15 2.5 1.25 "09.11. 20:45" $false

Set-BettingRow 17 16 "Italija 1" "Inter M." "Bologna" `
    1.4 5.2 7.6 1.1 3.1 1.2 "09.11. 20:45" $true

Set-BettingRow 18 17 "Italija 1" "H.Verona" "Juventus" `
    5 3.8 1.75 2.15 1.2 1.3 "09.11. 20:45" $true

Set-BettingRow 19 18 "Italija 1" "Lazio" "Monza" `
    1.65 4.2 5 1.2 2.3 1.25 "09.11. 20:45" $false

Set-BettingRow 20 19 "Italija 1" "Sampdoria" "Lecce" `
    2.6 3.1 2.95 1.4 1.5 1.4 "13.11. 20:45" $false

$ws.Range("L18").Select() | Out-Null
